$d = $word.ActiveDocument

$replacements = @(
    @("489×4=1956", "236×7=1652"),
    @("394×3=1182", "362×6=2172"),
    @("460×4=1840", "520×7=3640"),
    @("848×8=6784", "422×5=2110"),
    @("395×4=1580", "576×7=4032"),
    @("289×8=2312", "855×5=4275"),
    @("639×6=3834", "690×4=2760"),
    @("208×9=1872", "685×9=6165"),
    @("959×2=1918", "376×7=2632"),
    @("721×8=5768", "180×7=1260"),
    @("138×2=276", "300×9=2700"),
    @("659×5=3295", "375×8=3000"),
    @("644×8=5152", "604×3=1812"),
    @("468×7=3276", "214×7=1498"),
    @("970×2=1940", "765×7=5355"),
    @("126×3=378", "553×5=2765"),
    @("270×9=2430", "419×2=838"),
    @("564×4=2256", "802×4=3208"),
    @("326×2=652", "981×9=8829"),
    @("195×3=585", "707×2=1414"),
    @("246×2=492", "743×4=2972"),
    @("304×9=2736", "623×5=3115"),
    @("829×8=6632", "858×9=7722"),
    @("414×3=1242", "521×6=3126"),
    @("623×4=2492", "746×4=2984"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Done applying replacements."